$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.951.58"
$ws.Range("E2").Value = "  -3.39%  "
$ws.Range("D3").Value = "1.793.83"
$ws.Range("E3").Value = "  -3.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.71"
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4189"
$ws.Range("E7").Value = "  -3.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3554"
$ws.Range("E8").Value = "  -4.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07079"
$ws.Range("E9").Value = "  -4.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8436"
$ws.Range("E10").Value = "  -4.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.10"
$ws.Range("E11").Value = "  -5.21%  "
$ws.Range("D12").Value = "1.781.20"
$ws.Range("E12").Value = "  -2.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.278"
$ws.Range("E13").Value = "  -3.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.330"
$ws.Range("E14").Value = "  -4.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06757"
$ws.Range("E15").Value = "  -2.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.55"
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008685"
$ws.Range("E18").Value = "  -4.65%  "
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.97"
$ws.Range("E20").Value = "  -4.38%  "
$ws.Range("D21").Value = "26.951.03"
$ws.Range("E21").Value = "  -3.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.043"
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.93"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("D24").Value = "2.010.92"
$ws.Range("E24").Value = "  -3.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.932"
$ws.Range("E25").Value = "  -1.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.54"
$ws.Range("E26").Value = "  -1.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.05"
$ws.Range("E27").Value = "  -3.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.987"
$ws.Range("E28").Value = "  -6.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.86"
$ws.Range("E29").Value = "  -2.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.637"
$ws.Range("E30").Value = "  -12.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08922"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.859"
$ws.Range("E32").Value = "  -4.02%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7147"
$ws.Range("E33").Value = "  -9.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.280"
$ws.Range("E34").Value = "  -7.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.002"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.068"
$ws.Range("E36").Value = "  -9.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.074"
$ws.Range("E37").Value = "  -3.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01900"
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05091"
$ws.Range("E39").Value = "  -6.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.1621"
$ws.Range("E40").Value = "  -4.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4932"
$ws.Range("E41").Value = "  -5.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.573"
$ws.Range("E42").Value = "  -9.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.987"
$ws.Range("E43").Value = "  -11.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.997"
$ws.Range("E44").Value = "  -8.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.39"
$ws.Range("E45").Value = "  -2.59%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.20"
$ws.Range("E46").Value = "  -3.80%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("E48").Value = "  -4.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4491"
$ws.Range("E49").Value = "  -6.58%  "
$ws.Range("E50").Value = "  -5.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "61.81"
$ws.Range("E51").Value = "  -5.26%  "
